$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Fechas de la campaña para Constelación de botas 2022: 14-23 de mayo, 13-22 de junio, 12-21 de julio",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "2022 Fechas de la campaña para Constelación de botas: 14-23 de mayo, 13-22 de junio, 12-21 de julio",
    2
)
